$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.004958007358721677, 0.004861074973228985, 0.9903417954640102, 0.957446478056866, 0.03327798622130068)
    3  = @(0.003739810559255043, 0.003653747881690812, 0.9927148204291156, 0.9680131162517521, 0.01326005264751901)
    4  = @(0.003405827010953281, 0.003314319779483316, 0.9933654155001719, 0.9709827628938295, 0.01429345732811882)
    5  = @(0.003257129333634557, 0.003140833072028297, 0.9936550565119737, 0.9725002170295135, 0.01386102876658908)
    6  = @(0.003192107076096309, 0.003040422155846522, 0.9937817207758056, 0.9733784427813837, 0.01384516524969282)
    7  = @(0.003163942208685931, 0.00299182578862744,  0.9938365950173765, 0.9738030155393282, 0.01417268802553237)
    8  = @(0.0031532551802971,   0.002965733558919699, 0.993857414701213,  0.9740306548814219, 0.01401041383150893)
    9  = @(0.003149388712457963, 0.00295561973124334,  0.9938649480363753, 0.9741189114538299, 0.01384946254290247)
    10 = @(0.003147323221355497, 0.002954008572207624, 0.9938689686819906, 0.9741324929757319, 0.01373682671195912)
    11 = @(0.003145267823882309, 0.002955946713526951, 0.9938729814491797, 0.9741145774834474, 0.01378774163131053)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($col = 2; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowVals[$col - 2]
    }
}
